$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-21 and 23-51 (row 22 unchanged) need D/E updates.
# Force text format on D2:E51 so numeric-looking strings (e.g. "310.83")
# are stored as text, not auto-converted to numbers, then restore default style.
$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = '23.929.19'
$ws.Range("E2").Value = '  -1.86%  '
$ws.Range("D3").Value = '1.653.35'
$ws.Range("E3").Value = '  -0.70%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '310.83'
$ws.Range("E5").Value = '  -0.61%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("D7").Value = '0.3873'
$ws.Range("E7").Value = '  -1.91%  '
$ws.Range("D8").Value = '0.3819'
$ws.Range("E8").Value = '  -2.34%  '
$ws.Range("D9").Value = '51.89'
$ws.Range("E9").Value = '  -0.43%  '
$ws.Range("D10").Value = '1.352'
$ws.Range("E10").Value = '  -2.95%  '
$ws.Range("D11").Value = '1.000'
$ws.Range("E11").Value = '  -0.14%  '
$ws.Range("D12").Value = '0.08465'
$ws.Range("E12").Value = '  -1.16%  '
$ws.Range("D13").Value = '23.99'
$ws.Range("E13").Value = '  -1.35%  '
$ws.Range("D14").Value = '7.090'
$ws.Range("E14").Value = '  -2.64%  '
$ws.Range("D15").Value = '8.057'
$ws.Range("E15").Value = '  +1.50%  '
$ws.Range("D16").Value = '0.00001317'
$ws.Range("E16").Value = '  -1.61%  '
$ws.Range("D17").Value = '1.650.27'
$ws.Range("E17").Value = '  -0.74%  '
$ws.Range("D18").Value = '94.25'
$ws.Range("E18").Value = '  -0.75%  '
$ws.Range("D19").Value = '0.06997'
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("D20").Value = '19.69'
$ws.Range("E20").Value = '  -3.97%  '
$ws.Range("D21").Value = '6.961'
$ws.Range("E21").Value = '  -0.28%  '
$ws.Range("D23").Value = '13.80'
$ws.Range("E23").Value = '  +0.56%  '
$ws.Range("D24").Value = '23.925.88'
$ws.Range("E24").Value = '  -1.87%  '
$ws.Range("D25").Value = '2.435'
$ws.Range("E25").Value = '  +0.55%  '
$ws.Range("D26").Value = '2.988'
$ws.Range("E26").Value = '  -1.79%  '
$ws.Range("D27").Value = '22.11'
$ws.Range("E27").Value = '  -1.74%  '
$ws.Range("D28").Value = '153.78'
$ws.Range("E28").Value = '  -2.12%  '
$ws.Range("D29").Value = '5.423'
$ws.Range("E29").Value = '  -0.47%  '
$ws.Range("D30").Value = '137.92'
$ws.Range("E30").Value = '  -3.22%  '
$ws.Range("D31").Value = '7.878'
$ws.Range("E31").Value = '  -1.99%  '
$ws.Range("D32").Value = '2.504'
$ws.Range("E32").Value = '  -1.23%  '
$ws.Range("D33").Value = '1.833.34'
$ws.Range("E33").Value = '  -0.51%  '
$ws.Range("D34").Value = '1.022'
$ws.Range("E34").Value = '  -3.19%  '
$ws.Range("D35").Value = '0.08176'
$ws.Range("E35").Value = '  -0.49%  '
$ws.Range("D36").Value = '6.721'
$ws.Range("E36").Value = '  -2.26%  '
$ws.Range("D37").Value = '0.02925'
$ws.Range("E37").Value = '  -3.07%  '
$ws.Range("D38").Value = '10.83'
$ws.Range("E38").Value = '  -2.61%  '
$ws.Range("D39").Value = '0.2683'
$ws.Range("E39").Value = '  -2.58%  '
$ws.Range("D40").Value = '0.09126'
$ws.Range("E40").Value = '  -1.14%  '
$ws.Range("D41").Value = '0.7587'
$ws.Range("E41").Value = '  -1.47%  '
$ws.Range("D42").Value = '13.55'
$ws.Range("E42").Value = '  -1.49%  '
$ws.Range("D43").Value = '1.428'
$ws.Range("E43").Value = '  -1.36%  '
$ws.Range("D44").Value = '16.47'
$ws.Range("E44").Value = '  +0.07%  '
$ws.Range("D45").Value = '0.6951'
$ws.Range("E45").Value = '  -2.04%  '
$ws.Range("D46").Value = '2.466'
$ws.Range("E46").Value = '  -2.45%  '
$ws.Range("D47").Value = '4.100'
$ws.Range("E47").Value = '  -0.69%  '
$ws.Range("D48").Value = '1.000'
$ws.Range("E48").Value = '  +0.16%  '
$ws.Range("D49").Value = '0.08300'
$ws.Range("E49").Value = '  -1.41%  '
$ws.Range("D50").Value = '134.55'
$ws.Range("E50").Value = '  -1.65%  '
$ws.Range("D51").Value = '1.237'
$ws.Range("E51").Value = '  -2.30%  '

# Restore original (default) cell style so no stray formatting is introduced.
$rng.Style = "Normal"
